# Update BOM: only 1k resistor is needed (consolidate R1/R4 to 1k Ohm),
# consolidate Q3/Q4/Q5 transistor to Si2302DS, fix R13 to 680 Ohm, and
# add purchase-order hyperlinks/ranking numbers to the purchased-components
# section of the checklist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value (column C) updates ---
$ws.Range("C13").Value = "Si2302DS "
$ws.Range("C14").Value = "Si2302DS "
$ws.Range("C15").Value = "1k Ohm"
$ws.Range("C21").Value = "680 Ohm"

# --- DNP/order column (H) numeric purchase-order values ---
$ws.Range("H12").Value = 7
$ws.Range("H13").Value = 8
$ws.Range("H14").Value = 9
$ws.Range("H15").Value = 6
$ws.Range("H16").Value = 5
$ws.Range("H18").Value = 4
$ws.Range("H19").Value = 3
$ws.Range("H20").Value = 2
$ws.Range("H21").Value = 1

# --- Hyperlinks (column I); TextToDisplay keeps the visible cell text in
#     sync with the new link address (the sheet shows the raw URL as text). ---
$linkI3  = "https://makerselectronics.com/product/capacitor-470uf-35v-1016mm/?srsltid=AfmBOor09AGOe6uRfUPUSJ8gQ4mHe95vw-zHcd40OwpcpkuWNsINe4Bm"
$linkI13 = "https://uge-one.com/product/si2302ds-sot23-general-purpose-n-channel-mosfet-smd-transistor-sot-23/"
$linkI14 = "https://uge-one.com/product/si2302ds-sot23-general-purpose-n-channel-mosfet-smd-transistor-sot-23/"
$linkI15 = "https://uge-one.com/product/smd-chip-resistor-size-1206-1k-ohm-1/"
$linkI16 = "https://uge-one.com/product/smd-chip-resistor-size-1206-100r-ohm/"
$linkI18 = "https://uge-one.com/product/smd-chip-resistor-size-1206-1k-ohm-1/"
$linkI19 = "https://uge-one.com/product/smd-chip-resistor-size-1206-300r-ohm/"
$linkI20 = "https://uge-one.com/product/smd-chip-resistor-size-1206-220r-ohm/"
$linkI21 = "https://uge-one.com/product/smd-chip-resistor-size-1206-680r-ohm/"

$ws.Hyperlinks.Add($ws.Range("I3"), $linkI3, "", "", $linkI3) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I13"), $linkI13, "", "", $linkI13) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I14"), $linkI14, "", "", $linkI14) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I15"), $linkI15, "", "", $linkI15) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I16"), $linkI16, "", "", $linkI16) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I18"), $linkI18, "", "", $linkI18) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I19"), $linkI19, "", "", $linkI19) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I20"), $linkI20, "", "", $linkI20) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I21"), $linkI21, "", "", $linkI21) | Out-Null

# --- View: zoom to 85% and move selection to C14 ---
$ws.Range("C14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
